# Append a new row (2025-04-13 price) to each "Solar_Prices" sheet,
# mirroring the existing pattern of dated-row entries stored as text.

$wb = $excel.ActiveWorkbook

# Sheet name -> value to place in column B for the new 2025-04-13 row.
# (Same value as the prior day, 2025-04-12, per the source diff.)
$updates = [ordered]@{
    "N-Dense"                   = "40"
    "N-Type"                    = "41.5"
    "N-type Wafer"               = "1.25"
    "Cell Topcon 183mm"          = "0.303"
    "Module Topcon 183mm"        = "0.1"
    "Silver Rear_side"           = "5,192"
    "Silver Busbar front-side"   = "7,773"
    "Silver finger front-side"   = "7,823"
    "USD_CNY"                    = "7.3258"
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    $newRow = 43

    # Leading apostrophe forces these to be stored as text, matching the
    # existing inline-string cells used throughout these sheets.
    $ws.Range("A$newRow").Value = "'2025-04-13"
    $ws.Range("B$newRow").Value = "'" + $updates[$sheetName]
}
